# Refresh cryptos list: updated Price and Volume(1h) columns (and, where the
# underlying rank changed, the Coin name / Link too) from the latest scrape.
#
# A handful of Price values are bare decimals (e.g. "549.72") that Excel's
# Range.Value setter would otherwise auto-convert to a Number (and silently
# drop significant trailing zeros, e.g. "25.50" -> 25.5). Prefixing the
# literal with an apostrophe forces text entry exactly as typing '549.72 into
# the cell would in the Excel UI; resetting the Style back to Normal afterwards
# clears the quote-prefix formatting Excel applies so the cell keeps its
# original (unstyled) look while still holding a text value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.383.70'
$ws.Range("E2").Value = '  -0.56%  '
$ws.Range("D3").Value = '2.396.42'
$ws.Range("E3").Value = '  -3.68%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '''549.72'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.82%  '
$ws.Range("E6").Value = '  -3.04%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = '''0.539'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -10.38%  '
$ws.Range("D9").Value = '2.396.34'
$ws.Range("E9").Value = '  -3.64%  '
$ws.Range("E10").Value = '  -1.90%  '
$ws.Range("E11").Value = '  +0.23%  '
$ws.Range("D12").Value = '''5.31'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.25%  '
$ws.Range("E13").Value = '  -3.34%  '
$ws.Range("D14").Value = '''25.50'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.08%  '
$ws.Range("D15").Value = '2.829.78'
$ws.Range("E15").Value = '  -3.65%  '
$ws.Range("E16").Value = '  -1.38%  '
$ws.Range("D17").Value = '61.028.21'
$ws.Range("E17").Value = '  -1.02%  '
$ws.Range("D18").Value = '2.399.13'
$ws.Range("E18").Value = '  -3.73%  '
$ws.Range("D19").Value = '''10.82'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.33%  '
$ws.Range("E20").Value = '  -1.41%  '
$ws.Range("D21").Value = '''320.10'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.96%  '
$ws.Range("D22").Value = '''6.75'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.62%  '
$ws.Range("B23").Value = 'SuiNetwork'
$ws.Range("C23").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D23").Value = '''1.94'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +7.22%  '
$ws.Range("B24").Value = 'Dai'
$ws.Range("C24").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D24").Value = '''1.00'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.00%  '
$ws.Range("D25").Value = '''63.89'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.55%  '
$ws.Range("D26").Value = '''8.22'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +8.12%  '
$ws.Range("D27").Value = '''543.27'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.10%  '
$ws.Range("B28").Value = 'PEPE'
$ws.Range("C28").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D28").Value = '0.0₃0947'
$ws.Range("E28").Value = '  -4.72%  '
$ws.Range("B29").Value = 'Binance-PegBSC-USD'
$ws.Range("C29").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D29").Value = '''0.999'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.10%  '
$ws.Range("D30").Value = '2.514.83'
$ws.Range("E30").Value = '  -3.76%  '
$ws.Range("D31").Value = '''1.46'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.53%  '
$ws.Range("D32").Value = '''8.13'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.32%  '
$ws.Range("E33").Value = '  -3.36%  '
$ws.Range("E34").Value = '  -2.97%  '
$ws.Range("E35").Value = '  -0.52%  '
$ws.Range("D36").Value = '''1.00'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.02%  '
$ws.Range("E37").Value = '  -6.16%  '
$ws.Range("D38").Value = '''4.75'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.26%  '
$ws.Range("B39").Value = 'PolygonEcosystemToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D39").Value = '''0.378'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.52%  '
$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").Value = '''1.87'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +7.33%  '
$ws.Range("D41").Value = '''18.17'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.26%  '
$ws.Range("D42").Value = '''138.98'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -6.45%  '
$ws.Range("E43").Value = '  +0.03%  '
$ws.Range("D44").Value = '''40.21'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.59%  '
$ws.Range("D45").Value = '''2.23'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.93%  '
$ws.Range("D46").Value = '''142.15'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.44%  '
$ws.Range("D47").Value = '''3.64'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.00%  '
$ws.Range("D48").Value = '''20.40'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.70%  '
$ws.Range("D49").Value = '''0.0524'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.74%  '
$ws.Range("D50").Value = '''0.580'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.32%  '
$ws.Range("E51").Value = '  -1.13%  '
